# Deploying to gh-pages from @ w4bo/slides-markdown@5ffc620721c7f2e08da49fa862785a9bd553458f 🚀
#
# 1) Add a new "References" slide at the end of the deck (after the
#    current last slide), using the same "Title and Content" layout as
#    the other content slides.
# 2) Update the "Left column" text on slide 2 to credit the reference.

$p = $ppt.ActivePresentation

# --- 1. Append the new References slide -----------------------------------

$lastIndex = $p.Slides.Count
$layout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"
$newSlide = $p.Slides.AddSlide($lastIndex + 1, $layout)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "References"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange

$run1 = "Francia, Matteo, Enrico Gallinucci, and Matteo Golfarelli. 2024. " + [char]0x201C + "Colossal Trajectory Mining: A Unifying Approach to Mine Behavioral Mobility Patterns." + [char]0x201D + " "
$run2 = "Expert Syst. Appl."
$run3 = " 238 (Part E): 122055. "
$run4 = "https://doi.org/10.1016/J.ESWA.2023.122055"
$run5 = "."

$body.Text = $run1 + $run2 + $run3 + $run4 + $run5

$italicStart = $run1.Length + 1
$body.Characters($italicStart, $run2.Length).Font.Italic = $true

$linkStart = $run1.Length + $run2.Length + $run3.Length + 1
$linkRange = $body.Characters($linkStart, $run4.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = $run4

# --- 2. Update "Left column" text on slide 2 -------------------------------

$slide2 = $p.Slides.Item(2)
$slide2.Shapes.Item(2).TextFrame.TextRange.Text = "Left column (Francia, Gallinucci, and Golfarelli 2024)"
